$wb = $excel.ActiveWorkbook

# ALC row 64
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H64").Value = 3583.4666
$ws.Range("I64").Value = 3425.7778
$ws.Range("J64").Value = 3820
$ws.Range("K64").Value = 3425.7778
$ws.Range("L64").Value = 3820
$ws.Range("M64").Value = -3177.7778
$ws.Range("N64").Value = -4316

# ALC row 67
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H67").Value = 3583.4666
$ws.Range("I67").Value = 3425.7778
$ws.Range("J67").Value = 3820
$ws.Range("K67").Value = 3425.7778
$ws.Range("L67").Value = 3820
$ws.Range("M67").Value = -2567.7778
$ws.Range("N67").Value = -5536

# ALC row 138
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H138").Value = 5257.5728
$ws.Range("I138").Value = 3224.5386
$ws.Range("J138").Value = 6012.7
$ws.Range("K138").Value = 9673.6158
$ws.Range("L138").Value = 18038.1
$ws.Range("M138").Value = -4533.6158
$ws.Range("N138").Value = -28318.1

# ARM row 32
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 1416292.1
$ws.Range("I32").Value = 14027.448
$ws.Range("J32").Value = 21749130
$ws.Range("K32").Value = 14027.448
$ws.Range("L32").Value = 21749130
$ws.Range("M32").Value = -13740.448
$ws.Range("N32").Value = -21749704

# ARM row 61
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 2292.7273
$ws.Range("I61").Value = 2513
$ws.Range("J61").Value = 1705.3334
$ws.Range("K61").Value = 2513
$ws.Range("L61").Value = 1705.3334
$ws.Range("M61").Value = -2301
$ws.Range("N61").Value = -2129.3334

# ARM row 132
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H132").Value = 48569
$ws.Range("I132").Value = 2948.5625
$ws.Range("J132").Value = 170223.5
$ws.Range("K132").Value = 8845.6875
$ws.Range("L132").Value = 510670.5
$ws.Range("M132").Value = -6315.6875
$ws.Range("N132").Value = -515730.5

# ARM row 133
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H133").Value = 70558.53999999999
$ws.Range("I133").Value = 0
$ws.Range("J133").Value = 70558.53999999999
$ws.Range("K133").Value = 0
$ws.Range("L133").Value = 70558.53999999999
$ws.Range("N133").Value = -75618.53999999999

# ARM row 135
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H135").Value = 45631.89
$ws.Range("I135").Value = 0
$ws.Range("J135").Value = 45631.89
$ws.Range("K135").Value = 0
$ws.Range("L135").Value = 45631.89
$ws.Range("N135").Value = -55771.89

# ARM row 136
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H136").Value = 2292.7273
$ws.Range("I136").Value = 2513
$ws.Range("J136").Value = 1705.3334
$ws.Range("K136").Value = 7539
$ws.Range("L136").Value = 5116.0002
$ws.Range("M136").Value = -4989
$ws.Range("N136").Value = -10216.0002

# ARM row 137
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H137").Value = 100001
$ws.Range("I137").Value = 100001
$ws.Range("J137").Value = 100001
$ws.Range("K137").Value = 100001
$ws.Range("L137").Value = 100001
$ws.Range("M137").Value = -94901
$ws.Range("N137").Value = -110201

# ARM row 138
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H138").Value = 72000
$ws.Range("I138").Value = 0
$ws.Range("J138").Value = 72000
$ws.Range("K138").Value = 0
$ws.Range("L138").Value = 72000
$ws.Range("N138").Value = -82280

# ARM row 139
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H139").Value = 69000
$ws.Range("I139").Value = 0
$ws.Range("J139").Value = 69000
$ws.Range("K139").Value = 0
$ws.Range("L139").Value = 69000
$ws.Range("N139").Value = -79280

# ARM row 140
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H140").Value = 0
$ws.Range("I140").Value = 0
$ws.Range("J140").Value = 0
$ws.Range("K140").Value = 0
$ws.Range("L140").Value = ""
$ws.Range("N140").Value = 0

# ARM row 141
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H141").Value = 66000
$ws.Range("I141").Value = 0
$ws.Range("J141").Value = 66000
$ws.Range("K141").Value = 0
$ws.Range("L141").Value = 66000
$ws.Range("N141").Value = -76360

# BSM row 42
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H42").Value = 0
$ws.Range("I42").Value = 0
$ws.Range("J42").Value = 0
$ws.Range("K42").Value = 0
$ws.Range("L42").Value = ""
$ws.Range("N42").Value = 0

# BSM row 57
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H57").Value = 0
$ws.Range("I57").Value = 0
$ws.Range("J57").Value = 0
$ws.Range("K57").Value = 0
$ws.Range("L57").Value = ""
$ws.Range("N57").Value = 0

# BSM row 132
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H132").Value = 0
$ws.Range("I132").Value = 0
$ws.Range("J132").Value = 0
$ws.Range("K132").Value = 0
$ws.Range("L132").Value = ""
$ws.Range("N132").Value = 0

# BSM row 134
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 11211.357
$ws.Range("I134").Value = 1013.25
$ws.Range("J134").Value = 72400
$ws.Range("K134").Value = 3039.75
$ws.Range("L134").Value = 217200
$ws.Range("M134").Value = -504.75
$ws.Range("N134").Value = -222270

# BSM row 135
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H135").Value = 48963.332
$ws.Range("I135").Value = 0
$ws.Range("J135").Value = 48963.332
$ws.Range("K135").Value = 0
$ws.Range("L135").Value = 48963.332
$ws.Range("N135").Value = -59103.332

# BSM row 136
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H136").Value = 0
$ws.Range("I136").Value = 0
$ws.Range("J136").Value = 0
$ws.Range("K136").Value = 0
$ws.Range("L136").Value = ""
$ws.Range("N136").Value = 0

# BSM row 137
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H137").Value = 0
$ws.Range("I137").Value = 0
$ws.Range("J137").Value = 0
$ws.Range("K137").Value = 0
$ws.Range("L137").Value = ""
$ws.Range("N137").Value = 0

# BSM row 138
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H138").Value = 43170
$ws.Range("I138").Value = 0
$ws.Range("J138").Value = 43170
$ws.Range("K138").Value = 0
$ws.Range("L138").Value = 43170
$ws.Range("N138").Value = -53450

# BSM row 139
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H139").Value = 0
$ws.Range("I139").Value = 0
$ws.Range("J139").Value = 0
$ws.Range("K139").Value = 0
$ws.Range("L139").Value = ""
$ws.Range("N139").Value = 0

# BSM row 140
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H140").Value = 86795
$ws.Range("I140").Value = 0
$ws.Range("J140").Value = 86795
$ws.Range("K140").Value = 0
$ws.Range("L140").Value = 86795
$ws.Range("N140").Value = -97155

# BSM row 141
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H141").Value = 116593.336
$ws.Range("I141").Value = 0
$ws.Range("J141").Value = 116593.336
$ws.Range("K141").Value = 0
$ws.Range("L141").Value = 116593.336
$ws.Range("M141").Value = ""
$ws.Range("N141").Value = -126953.336

# CRP row 31
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 7979.411
$ws.Range("I31").Value = 2298.8833
$ws.Range("J31").Value = 19340.467
$ws.Range("K31").Value = 2298.8833
$ws.Range("L31").Value = 19340.467
$ws.Range("M31").Value = -2003.8833
$ws.Range("N31").Value = -19930.467

# CRP row 34
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H34").Value = 7979.411
$ws.Range("I34").Value = 2298.8833
$ws.Range("J34").Value = 19340.467
$ws.Range("K34").Value = 2298.8833
$ws.Range("L34").Value = 19340.467
$ws.Range("M34").Value = -2096.8833
$ws.Range("N34").Value = -19744.467

# CRP row 54
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H54").Value = 8000
$ws.Range("I54").Value = 0
$ws.Range("J54").Value = 8000
$ws.Range("K54").Value = 0
$ws.Range("L54").Value = 8000
$ws.Range("N54").Value = -9316

# CRP row 137
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H137").Value = 48556.668
$ws.Range("I137").Value = 0
$ws.Range("J137").Value = 48556.668
$ws.Range("K137").Value = 0
$ws.Range("L137").Value = 48556.668
$ws.Range("N137").Value = -58756.668

# CRP row 140
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H140").Value = 89350
$ws.Range("I140").Value = 0
$ws.Range("J140").Value = 89350
$ws.Range("K140").Value = 0
$ws.Range("L140").Value = 89350
$ws.Range("N140").Value = -99710

# CUL row 5
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 1286
$ws.Range("I5").Value = 995
$ws.Range("J5").Value = 2450
$ws.Range("K5").Value = 2985
$ws.Range("L5").Value = 7350
$ws.Range("M5").Value = -2873
$ws.Range("N5").Value = -7574

# CUL row 110
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H110").Value = 5113.3887
$ws.Range("I110").Value = 1286.375
$ws.Range("J110").Value = 8175
$ws.Range("K110").Value = 3859.125
$ws.Range("L110").Value = 24525
$ws.Range("M110").Value = 230.875
$ws.Range("N110").Value = -32705

# CUL row 135
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H135").Value = 1286
$ws.Range("I135").Value = 995
$ws.Range("J135").Value = 2450
$ws.Range("K135").Value = 8955
$ws.Range("L135").Value = 22050
$ws.Range("M135").Value = -6420
$ws.Range("N135").Value = -27120

# LTW row 40
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 15387359
$ws.Range("I40").Value = 2648
$ws.Range("J40").Value = 40002896
$ws.Range("K40").Value = 2648
$ws.Range("L40").Value = 40002896
$ws.Range("M40").Value = -2512
$ws.Range("N40").Value = -40003168

# LTW row 132
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H132").Value = 5819.294
$ws.Range("I132").Value = 6384
$ws.Range("J132").Value = 4464
$ws.Range("K132").Value = 19152
$ws.Range("L132").Value = 13392
$ws.Range("M132").Value = -16622
$ws.Range("N132").Value = -18452

# LTW row 136
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H136").Value = 2474.197
$ws.Range("I136").Value = 1529.6976
$ws.Range("J136").Value = 4240
$ws.Range("K136").Value = 4589.0928
$ws.Range("L136").Value = 12720
$ws.Range("M136").Value = -2039.0928
$ws.Range("N136").Value = -17820

# LTW row 139
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H139").Value = 79550
$ws.Range("I139").Value = 0
$ws.Range("J139").Value = 79550
$ws.Range("K139").Value = 0
$ws.Range("L139").Value = 79550
$ws.Range("N139").Value = -89830

# WVR row 136
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H136").Value = 1238.3112
$ws.Range("I136").Value = 776.0357
$ws.Range("J136").Value = 1999.7059
$ws.Range("K136").Value = 2328.1071
$ws.Range("L136").Value = 5999.1177
$ws.Range("M136").Value = 221.8928999999998
$ws.Range("N136").Value = -11099.1177
